$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.295.70"
$ws.Range("E2").Value = "'  -3.21%  "
$ws.Range("D3").Value = "'1.740.99"
$ws.Range("E3").Value = "'  -3.29%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'321.08"
$ws.Range("E5").Value = "'  -4.26%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("D7").Value = "'0.4200"
$ws.Range("E7").Value = "'  -7.79%  "
$ws.Range("D8").Value = "'0.3569"
$ws.Range("E8").Value = "'  -3.86%  "
$ws.Range("D9").Value = "'45.02"
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("D10").Value = "'0.07367"
$ws.Range("E10").Value = "'  -2.87%  "
$ws.Range("D11").Value = "'1.109"
$ws.Range("E11").Value = "'  -3.38%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "'  +0.13%  "
$ws.Range("D13").Value = "'21.40"
$ws.Range("E13").Value = "'  -4.00%  "
$ws.Range("D14").Value = "'6.047"
$ws.Range("E14").Value = "'  -4.33%  "
$ws.Range("D15").Value = "'7.157"
$ws.Range("E15").Value = "'  -4.33%  "
$ws.Range("D16").Value = "'1.737.45"
$ws.Range("E16").Value = "'  -3.26%  "
$ws.Range("D17").Value = "'0.00001064"
$ws.Range("E17").Value = "'  -3.02%  "
$ws.Range("D18").Value = "'86.12"
$ws.Range("E18").Value = "'  +5.27%  "
$ws.Range("D19").Value = "'0.06010"
$ws.Range("E19").Value = "'  -10.42%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "'  +0.08%  "
$ws.Range("D21").Value = "'16.74"
$ws.Range("E21").Value = "'  -4.18%  "
$ws.Range("D22").Value = "'6.062"
$ws.Range("E22").Value = "'  -5.19%  "
$ws.Range("D23").Value = "'0.5240"
$ws.Range("E23").Value = "'  -7.04%  "
$ws.Range("D24").Value = "'27.341.70"
$ws.Range("E24").Value = "'  -2.94%  "
$ws.Range("D25").Value = "'11.33"
$ws.Range("E25").Value = "'  -4.44%  "
$ws.Range("D26").Value = "'2.355"
$ws.Range("E26").Value = "'  -2.24%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.33"
$ws.Range("E27").Value = "'  -1.37%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'152.73"
$ws.Range("E28").Value = "'  +0.52%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.373"
$ws.Range("E29").Value = "'  +0.26%  "
$ws.Range("D30").Value = "'1.938.56"
$ws.Range("E30").Value = "'  -3.10%  "
$ws.Range("D31").Value = "'125.68"
$ws.Range("E31").Value = "'  -5.75%  "
$ws.Range("D32").Value = "'1.165"
$ws.Range("E32").Value = "'  -6.71%  "
$ws.Range("D33").Value = "'5.664"
$ws.Range("E33").Value = "'  -3.02%  "
$ws.Range("D34").Value = "'0.09045"
$ws.Range("E34").Value = "'  -5.07%  "
$ws.Range("D35").Value = "'3.601"
$ws.Range("E35").Value = "'  -10.81%  "
$ws.Range("D36").Value = "'12.61"
$ws.Range("E36").Value = "'  +5.18%  "
$ws.Range("D37").Value = "'0.2132"
$ws.Range("E37").Value = "'  -3.72%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02261"
$ws.Range("E38").Value = "'  -3.90%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.030"
$ws.Range("E39").Value = "'  -3.79%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06018"
$ws.Range("E40").Value = "'  -5.14%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6337"
$ws.Range("E41").Value = "'  -4.49%  "
$ws.Range("D42").Value = "'1.185"
$ws.Range("E42").Value = "'  -3.90%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.428"
$ws.Range("E43").Value = "'  -5.44%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "'  +0.08%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.903"
$ws.Range("E45").Value = "'  -2.22%  "
$ws.Range("D46").Value = "'13.72"
$ws.Range("E46").Value = "'  -3.25%  "
$ws.Range("D47").Value = "'3.700"
$ws.Range("E47").Value = "'  -3.18%  "
$ws.Range("D48").Value = "'0.5776"
$ws.Range("E48").Value = "'  -4.94%  "
$ws.Range("D49").Value = "'124.84"
$ws.Range("E49").Value = "'  -4.35%  "
$ws.Range("D50").Value = "'1.941"
$ws.Range("E50").Value = "'  -5.01%  "
$ws.Range("D51").Value = "'0.06814"
$ws.Range("E51").Value = "'  -4.42%  "
